# TC02_Canine_Filter_Breed-AmerStaffd_Neo4jData.xlsx
# "updated all canine test cases with function changes"
#
# Two sheets are touched by this edit:
#
#   1. CaseDetailStat - the one-row header (File Name / File Type /
#      Association / Description / Format / Size) is removed, leaving the
#      header row blank while the data row underneath is untouched.
#
#   2. CaseDetailStat_Message - the third Neo4j connection/message block
#      (rows 21-30) now starts with a "Cypher query should not be an empty
#      string" warning (a new row, pushing the rest of the block down by
#      one), and the Cypher query text that used to be logged for that
#      block is now blank because the query string itself ended up empty.

$wb = $excel.ActiveWorkbook

# --- 1. CaseDetailStat: blank out the header row (row 1) ---
$wsStat = $wb.Worksheets.Item("CaseDetailStat")
$wsStat.Range("A1:F1").ClearContents()

# --- 2. CaseDetailStat_Message: new warning row + blanked query text ---
$wsMsg = $wb.Worksheets.Item("CaseDetailStat_Message")

# Shift the third message block (old rows 21-30) down by inserting a new
# row 21 ahead of it.
$wsMsg.Rows.Item(21).Insert()
$wsMsg.Range("A21").Value = "Cypher query should not be an empty string"

# The Cypher query text logged for that block (now at row 29, previously
# row 28) becomes an empty string.
$wsMsg.Range("A29").Formula = '=""'
